# Reddit X-ray CWL Random Distribution -- September CWL update
#
# - Columns A and C held placeholder cells (a quote-prefixed empty string,
#   used only so the "eligible" / "ineligible (missed hits)" columns had
#   something to apply error-handling formatting to); they become plain
#   blank cells now that the roster is driven entirely from column B.
# - Column B's roster is replaced with the new (longer) list of names for
#   "Ineligible (Already received)".
# - The used range grows from row 25 down to row 27 to fit the new roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Lord Zameow",
    "Bad Apple",
    "ANBU",
    "Seven Thunders",
    "Arcohol",
    "Annayake",
    "Protips",
    "Satan",
    "tre sedicesimi",
    "Nitin 4.0"
)

$firstRow = 2
$lastRow = 27
$xlPasteFormats = -4122

# A clean, un-prefixed style donor: the header cells use the plain
# "general alignment" style with no stale quote-prefix flag. Copying its
# format onto A/B/C re-stamps every data-row cell with that clean style
# (this also materializes rows 26-27, which don't exist yet).
$styleDonor = $ws.Range("B1")
$styleDonor.Copy()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)

    # Match the existing custom row height used by the rest of the table.
    $ws.Rows($r).RowHeight = 18.75

    $ws.Cells.Item($r, 1).Value = $null
    $ws.Cells.Item($r, 3).Value = $null

    $idx = $r - $firstRow
    if ($idx -lt $names.Length) {
        $ws.Cells.Item($r, 2).Value = $names[$idx]
    } else {
        $ws.Cells.Item($r, 2).Value = $null
    }
}

$excel.CutCopyMode = $false
